$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values scraped from source for the cryptos list.
# Numeric-looking text values (e.g. "0.999", "7.00") are written with a
# leading apostrophe so Excel keeps them as text instead of converting
# them to numbers (which would drop formatting like trailing zeros).

$ws.Range("D2").Value = '65.262.84'
$ws.Range("D3").Value = '3.407.64'
$ws.Range("E3").Value = '  -3.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''582.31'
$ws.Range("E5").Value = '  -2.84%  '
$ws.Range("D6").Value = '''136.33'
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.406.22'
$ws.Range("E8").Value = '  -3.41%  '
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  -7.57%  '
$ws.Range("E11").Value = '  -10.63%  '
$ws.Range("D12").Value = '''0.373'
$ws.Range("E12").Value = '  -7.81%  '
$ws.Range("D13").Value = '3.984.23'
$ws.Range("E13").Value = '  -3.51%  '
$ws.Range("D14").Value = '''0.0000178'
$ws.Range("E14").Value = '  -10.96%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '''0.115'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.410.71'
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("D17").Value = '65.173.04'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '''25.81'
$ws.Range("E18").Value = '  -9.91%  '
$ws.Range("D19").Value = '''9.72'
$ws.Range("E19").Value = '  -10.87%  '
$ws.Range("E20").Value = '  -5.75%  '
$ws.Range("D21").Value = '''13.50'
$ws.Range("E21").Value = '  -5.88%  '
$ws.Range("D22").Value = '''383.94'
$ws.Range("E22").Value = '  -7.28%  '
$ws.Range("D23").Value = '''0.553'
$ws.Range("E23").Value = '  -7.50%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '''72.40'
$ws.Range("E25").Value = '  -6.32%  '
$ws.Range("D26").Value = '3.542.49'
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("E27").Value = '  -10.24%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '  -10.15%  '
$ws.Range("D30").Value = '''7.00'
$ws.Range("E30").Value = '  -10.18%  '
$ws.Range("D31").Value = '''8.01'
$ws.Range("E31").Value = '  -10.15%  '
$ws.Range("D32").Value = '3.411.03'
$ws.Range("E32").Value = '  -3.28%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -7.63%  '
$ws.Range("D35").Value = '''22.69'
$ws.Range("E35").Value = '  -6.75%  '
$ws.Range("D36").Value = '''170.51'
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").Value = '''6.72'
$ws.Range("E37").Value = '  -10.68%  '
$ws.Range("D38").Value = '''1.45'
$ws.Range("E38").Value = '  -8.26%  '
$ws.Range("E39").Value = '  -13.58%  '
$ws.Range("D40").Value = '''4.69'
$ws.Range("E40").Value = '  -10.82%  '
$ws.Range("D41").Value = '''0.0755'
$ws.Range("E41").Value = '  -7.80%  '
$ws.Range("E42").Value = '  -5.56%  '
$ws.Range("D43").Value = '''43.41'
$ws.Range("E43").Value = '  -4.25%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '''4.37'
$ws.Range("E45").Value = '  -13.72%  '
$ws.Range("E46").Value = '  -10.83%  '
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").Value = '''22.12'
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("D49").Value = '''6.48'
$ws.Range("E49").Value = '  -8.26%  '
$ws.Range("D50").Value = '''2.03'
$ws.Range("E50").Value = '  -15.79%  '
$ws.Range("D51").Value = '2.169.04'
$ws.Range("E51").Value = '  -7.81%  '
